# "Add files via upload" / "reading week"
#
# 1. Highlight the header row (A1:E1) with a solid blue fill.
# 2. Turn on AutoFilter over the table (A1:E4) — this also registers the
#    hidden workbook-level _xlnm._FilterDatabase defined name that Excel
#    creates whenever AutoFilter is applied.
# 3. Append three more rows of data (4, 5, 6) in column A below the table.
# 4. Tidy up the sheet view: select the header row and scroll back to the
#    top-left corner (the sheet had previously been scrolled/selected to
#    B1/E7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row fill -----------------------------------------------
# RGB(0,112,192) == 0x0070C0 -> OLE/COM colors are BGR, so 0x00C07000 = 12611584
$ws.Range("A1:E1").Interior.Color = 12611584

# --- 2. AutoFilter + its hidden defined name ---------------------------
# Apply the filter before appending the extra rows below, so the
# autoFilter/_FilterDatabase range stays pinned to the original table
# (A1:E4) instead of growing to include the newly appended rows.
$ws.Range("A1:E4").AutoFilter()

$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Лист1!`$A`$1:`$E`$4")
$filterName.Visible = $false

# --- 3. New rows of data -----------------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

# --- 4. Selection / view ------------------------------------------------
$ws.Range("A1:E1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
